$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 408, pushing existing rows 408.. down by one
$ws.Rows.Item(408).Insert()

# Populate the newly inserted row 408 with the new weekly price record
$ws.Cells.Item(408, 1).Value2  = 11
$ws.Cells.Item(408, 2).Value2  = "Vega Monumental Concepción"
$ws.Cells.Item(408, 3).Value2  = "Bíobío"
$ws.Cells.Item(408, 4).Value2  = 45211
$ws.Cells.Item(408, 5).Value2  = 8
$ws.Cells.Item(408, 6).Value2  = 100112009
$ws.Cells.Item(408, 7).Value2  = "Acelga"
$ws.Cells.Item(408, 8).Value2  = "Sin especificar"
$ws.Cells.Item(408, 9).Value2  = "Primera"
$ws.Cells.Item(408, 10).Value2 = 180
$ws.Cells.Item(408, 11).Value2 = 600
$ws.Cells.Item(408, 12).Value2 = 650
$ws.Cells.Item(408, 13).Value2 = 622
$ws.Cells.Item(408, 14).Value2 = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(408, 15).Value2 = "Región de Ñuble"
$ws.Cells.Item(408, 16).Value2 = 622
$ws.Cells.Item(408, 17).Value2 = 1
$ws.Cells.Item(408, 18).Value2 = "Hortaliza"
